$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$new1 = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.9 = 41683.58 pesos`n✅ 41683.58 pesos = 9.87 = 973.38 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $new1

# --- tasas: update the N10/O10/N12/O12 rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 100.999
$ws2.Range("O10").Value = 4210

$ws2.Range("N12").Value = 4224.99
$ws2.Range("O12").Value = 98.66
